# Applies the cryptos-list price/volume refresh described in the commit
# ("Updated cryptos list ... with GitHub Actions"): column D (Price) and
# column E (Volume(1h)) text values are refreshed for the affected rows.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($cell, $text) {
    # Preserve the pre-existing cell style while forcing the assigned
    # value to be stored as text (these columns hold inline strings like
    # '566.92' or '2.540.68' that Excel would otherwise coerce to numbers).
    $origStyle = $cell.Style
    $cell.NumberFormat = "@"
    $cell.Value = $text
    $cell.Style = $origStyle
}

Set-TextValue $ws.Range("D2") "62.744.16"
$ws.Range("E2").Value = "  -0.71%  "
Set-TextValue $ws.Range("D3") "2.540.68"
$ws.Range("E3").Value = "  +3.15%  "
$ws.Range("E4").Value = "  +0.02%  "
Set-TextValue $ws.Range("D5") "566.92"
$ws.Range("E5").Value = "  +0.04%  "
Set-TextValue $ws.Range("D6") "146.55"
$ws.Range("E6").Value = "  +2.22%  "
$ws.Range("E7").Value = "  +0.05%  "
$ws.Range("E8").Value = "  -1.11%  "
Set-TextValue $ws.Range("D9") "2.539.52"
$ws.Range("E9").Value = "  +3.12%  "
$ws.Range("E10").Value = "  -1.33%  "
$ws.Range("E11").Value = "  -2.11%  "
$ws.Range("E13").Value = "  -0.39%  "
Set-TextValue $ws.Range("D14") "27.18"
$ws.Range("E14").Value = "  +2.47%  "
Set-TextValue $ws.Range("D15") "2.996.12"
$ws.Range("E15").Value = "  +3.18%  "
Set-TextValue $ws.Range("D16") "62.779.90"
$ws.Range("E16").Value = "  -0.52%  "
$ws.Range("E17").Value = "  -0.61%  "
Set-TextValue $ws.Range("D18") "2.542.42"
$ws.Range("E18").Value = "  +3.19%  "
$ws.Range("E19").Value = "  +1.89%  "
Set-TextValue $ws.Range("D20") "335.27"
$ws.Range("E20").Value = "  -1.91%  "
$ws.Range("E21").Value = "  -0.53%  "
Set-TextValue $ws.Range("D22") "6.73"
$ws.Range("E22").Value = "  -1.12%  "
$ws.Range("E23").Value = "  +0.12%  "
Set-TextValue $ws.Range("D24") "65.25"
$ws.Range("E24").Value = "  -0.55%  "
$ws.Range("E25").Value = "  -3.08%  "
$ws.Range("E26").Value = "  +6.98%  "
$ws.Range("E27").Value = "  +11.68%  "
Set-TextValue $ws.Range("D28") "1.00"
$ws.Range("E28").Value = "  +0.16%  "
Set-TextValue $ws.Range("D29") "8.37"
$ws.Range("E29").Value = "  +2.82%  "
Set-TextValue $ws.Range("D30") "7.26"
$ws.Range("E30").Value = "  +5.96%  "
Set-TextValue $ws.Range("D31") "0.0₃0807"
$ws.Range("E31").Value = "  -0.70%  "
$ws.Range("E32").Value = "  -0.19%  "
Set-TextValue $ws.Range("D33") "177.16"
$ws.Range("E33").Value = "  +1.21%  "
$ws.Range("E34").Value = "  +3.63%  "
Set-TextValue $ws.Range("D35") "409.12"
$ws.Range("E35").Value = "  +10.69%  "
Set-TextValue $ws.Range("D36") "0.399"
$ws.Range("E36").Value = "  -0.22%  "
Set-TextValue $ws.Range("D37") "18.95"
$ws.Range("E37").Value = "  +0.11%  "
$ws.Range("E38").Value = "  -0.02%  "
Set-TextValue $ws.Range("D39") "4.35"
$ws.Range("E39").Value = "  -2.29%  "
$ws.Range("E40").Value = "  +3.03%  "
Set-TextValue $ws.Range("D41") "1.00"
$ws.Range("E41").Value = "  -0.05%  "
Set-TextValue $ws.Range("D42") "39.06"
$ws.Range("E42").Value = "  -3.34%  "
Set-TextValue $ws.Range("D43") "152.92"
$ws.Range("E43").Value = "  +0.68%  "
$ws.Range("E44").Value = "  +0.55%  "
Set-TextValue $ws.Range("D45") "20.76"
$ws.Range("E45").Value = "  +0.93%  "
$ws.Range("E46").Value = "  +0.92%  "
$ws.Range("E47").Value = "  -0.51%  "
Set-TextValue $ws.Range("D48") "0.0517"
$ws.Range("E48").Value = "  -0.57%  "
$ws.Range("E49").Value = "  +4.06%  "
Set-TextValue $ws.Range("D50") "18.19"
$ws.Range("E50").Value = "  +1.02%  "
$ws.Range("E51").Value = "  +0.09%  "
